$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.031.33"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.34%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.303.23"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.77%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.58"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.38"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.10%  "

$ws.Range("E7").Value = "  +2.36%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  -1.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.37"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.23%  "

$ws.Range("E11").Value = "  -0.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "17.97"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.36%  "

$ws.Range("E13").Value = "  +0.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.81"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.663.00"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.316.45"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.784"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.92%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.983.12"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.66"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0910"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.12"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.41"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.06"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.93%  "

$ws.Range("E24").Value = "  -1.08%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.44"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.75%  "

$ws.Range("E27").Value = "  -0.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.18"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.94"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.12%  "

$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("E31").Value = "  -0.88%  "

$ws.Range("E32").Value = "  -3.27%  "

$ws.Range("E33").Value = "  +0.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.77"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.38%  "

$ws.Range("E35").Value = "  -3.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.75"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.68%  "

$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0690"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.07%  "

$ws.Range("E39").Value = "  -1.51%  "

$ws.Range("E40").Value = "  -1.56%  "

$ws.Range("E41").Value = "  +0.87%  "

$ws.Range("E42").Value = "  +0.93%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.998.63"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.26%  "

$ws.Range("E44").Value = "  -1.61%  "

$ws.Range("E45").Value = "  -2.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.23"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.99%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.53"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.62%  "

$ws.Range("E48").Value = "  -3.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.64"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.530.06"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.99"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.87%  "

